$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (C), ADTV (D) and Vola (E) values per row
$updates = @(
    @{ Row = 2; C = 34.84;              D = 691737;  E = 0.200459709970345 },
    @{ Row = 3; C = 73.66;              D = 89881;   E = 0.1563302702407023 },
    @{ Row = 4; C = 1296.59;            D = 18912;   E = 0.4882997222866197 },
    @{ Row = 5; C = 103.13;             D = 64367;   E = 0.2025903436990393 },
    @{ Row = 6; C = 86.34;              D = 3444314; E = 0.1802370915375292 },
    @{ Row = 7; C = 480.66;             D = 364870;  E = 0.1541333210277314 },
    @{ Row = 8; C = 94.47;              D = 3949377; E = 0.1221163772216723 },
    @{ Row = 9; C = 227.87;             D = 1609479; E = 0.2033037717085315 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
